$d = $word.ActiveDocument
$d.Content.Find.Execute("predict feature crime", $true, $false, $false, $false, $false,
                         $true, 1, $false, "predict future crime", 2)
